$d = $word.ActiveDocument

# Remove the trailing "*" from step numbers 4, 5 and 6 in the test plan
# table (step 7 keeps its "*" - it is still a manual/output step).
$d.Content.Find.Execute("4*", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "4", 2)
$d.Content.Find.Execute("5*", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "5", 2)
$d.Content.Find.Execute("6*", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "6", 2)
